$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark execution mode for existing test cases as "Manual"
$ws.Range("E2").Value = "Manual"
$ws.Range("E3").Value = "Manual"

# Update the active selection to E4
$ws.Range("E4").Select()
